$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.70"
$ws.Range("E2").Value = "'6.65%"
$ws.Range("D3").Value = "'31.87"
$ws.Range("E3").Value = "'8.65%"
$ws.Range("D4").Value = "'5.268"
$ws.Range("E4").Value = "'3.77%"
$ws.Range("D5").Value = "'0.07495"
$ws.Range("E5").Value = "'11.57%"
$ws.Range("D6").Value = "'7.836"
$ws.Range("E6").Value = "'7.18%"
$ws.Range("D7").Value = "'3.755"
$ws.Range("E7").Value = "'9.10%"
$ws.Range("D8").Value = "'1.470"
$ws.Range("E8").Value = "'5.50%"
$ws.Range("D9").Value = "'0.9149"
$ws.Range("E9").Value = "'2.47%"
$ws.Range("D10").Value = "'0.01741"
$ws.Range("E10").Value = "'2,586.77%"
$ws.Range("D11").Value = "'0.1702"
$ws.Range("E11").Value = "'8.02%"
$ws.Range("D12").Value = "'0.07791"
$ws.Range("E12").Value = "'9.27%"
$ws.Range("D13").Value = "'0.08052"
$ws.Range("E13").Value = "'5.71%"
$ws.Range("D14").Value = "'0.02989"
$ws.Range("E14").Value = "'2.29%"
$ws.Range("D15").Value = "'0.09898"
$ws.Range("E15").Value = "'10.22%"
$ws.Range("D16").Value = "'0.001491"
$ws.Range("E16").Value = "'-6.08%"
$ws.Range("D17").Value = "'0.04560"
$ws.Range("E17").Value = "'1.78%"
$ws.Range("D18").Value = "'0.006133"
$ws.Range("E18").Value = "'-1.28%"
$ws.Range("D19").Value = "'3.496"
$ws.Range("E19").Value = "'1.34%"
$ws.Range("D20").Value = "'2.234"
$ws.Range("E20").Value = "'0.18%"
$ws.Range("D21").Value = "'0.3302"
$ws.Range("E21").Value = "'2.12%"
$ws.Range("D22").Value = "'0.1342"
$ws.Range("E22").Value = "'1.83%"
$ws.Range("D23").Value = "'4.482"
$ws.Range("E23").Value = "'12.14%"
$ws.Range("D24").Value = "'0.1618"
$ws.Range("E24").Value = "'3.78%"
$ws.Range("E25").Value = "'0.93%"
$ws.Range("D26").Value = "'0.004438"
$ws.Range("E26").Value = "'1.44%"
$ws.Range("D27").Value = "'0.0001397"
$ws.Range("E27").Value = "'19.48%"
$ws.Range("D28").Value = "'0.0001738"
$ws.Range("E28").Value = "'7.49%"
$ws.Range("D40").Value = "'0.04518"
$ws.Range("E40").Value = "'6.46%"
$ws.Range("D41").Value = "'0.007213"
$ws.Range("E41").Value = "'5.73%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'8.80%"
$ws.Range("D43").Value = "'0.002385"
$ws.Range("E43").Value = "'9.98%"
$ws.Range("D44").Value = "'0.01395"
$ws.Range("E44").Value = "'21.95%"
$ws.Range("D45").Value = "'0.00006209"
$ws.Range("E45").Value = "'8.00%"
$ws.Range("E46").Value = "'-62.92%"
$ws.Range("E47").Value = "'-13.55%"
